$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F (dSF)
$changes = @{
    2  = -2
    9  = -4
    10 = 2
    11 = 4
    14 = -2
    15 = -6
    17 = -1
    21 = 4
    30 = 6
    32 = -3
    35 = 4
    37 = 0
    38 = -3
    39 = -4
    42 = -6
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
